$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 216; existing rows 216-224 shift down to 219-227.
$ws.Range("A216:R218").EntireRow.Insert()

# New row 216: Choclero / Primera, Región de O'Higgins
$ws.Cells.Item(216,1).Value = 4
$ws.Cells.Item(216,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216,3).Value = "Los Lagos"
$ws.Cells.Item(216,4).Value = 44610
$ws.Cells.Item(216,5).Value = 10
$ws.Cells.Item(216,6).Value = 100112024
$ws.Cells.Item(216,7).Value = "Choclo"
$ws.Cells.Item(216,8).Value = "Choclero"
$ws.Cells.Item(216,9).Value = "Primera"
$ws.Cells.Item(216,10).Value = 7500
$ws.Cells.Item(216,11).Value = 420
$ws.Cells.Item(216,12).Value = 420
$ws.Cells.Item(216,13).Value = 420
$ws.Cells.Item(216,14).Value = "`$/unidad"
$ws.Cells.Item(216,15).Value = "Región de O'Higgins"
$ws.Cells.Item(216,16).Value = 420
$ws.Cells.Item(216,17).Value = 1
$ws.Cells.Item(216,18).Value = "Hortaliza"

# New row 217: Choclero / Segunda, Región de O'Higgins
$ws.Cells.Item(217,1).Value = 4
$ws.Cells.Item(217,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(217,3).Value = "Los Lagos"
$ws.Cells.Item(217,4).Value = 44610
$ws.Cells.Item(217,5).Value = 10
$ws.Cells.Item(217,6).Value = 100112024
$ws.Cells.Item(217,7).Value = "Choclo"
$ws.Cells.Item(217,8).Value = "Choclero"
$ws.Cells.Item(217,9).Value = "Segunda"
$ws.Cells.Item(217,10).Value = 7500
$ws.Cells.Item(217,11).Value = 270
$ws.Cells.Item(217,12).Value = 270
$ws.Cells.Item(217,13).Value = 270
$ws.Cells.Item(217,14).Value = "`$/unidad"
$ws.Cells.Item(217,15).Value = "Región de O'Higgins"
$ws.Cells.Item(217,16).Value = 270
$ws.Cells.Item(217,17).Value = 1
$ws.Cells.Item(217,18).Value = "Hortaliza"

# New row 218: Dulce o Americano / Primera, Región de O'Higgins
$ws.Cells.Item(218,1).Value = 4
$ws.Cells.Item(218,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(218,3).Value = "Los Lagos"
$ws.Cells.Item(218,4).Value = 44610
$ws.Cells.Item(218,5).Value = 10
$ws.Cells.Item(218,6).Value = 100112024
$ws.Cells.Item(218,7).Value = "Choclo"
$ws.Cells.Item(218,8).Value = "Dulce o Americano"
$ws.Cells.Item(218,9).Value = "Primera"
$ws.Cells.Item(218,10).Value = 25000
$ws.Cells.Item(218,11).Value = 220
$ws.Cells.Item(218,12).Value = 250
$ws.Cells.Item(218,13).Value = 232
$ws.Cells.Item(218,14).Value = "`$/unidad"
$ws.Cells.Item(218,15).Value = "Región de O'Higgins"
$ws.Cells.Item(218,16).Value = 232
$ws.Cells.Item(218,17).Value = 1
$ws.Cells.Item(218,18).Value = "Hortaliza"
